$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.566.68"
$ws.Range("E2").Value = "  -2.56%  "
$ws.Range("D3").Value = "3.429.60"
$ws.Range("E3").Value = "  -5.34%  "
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").Value = "'571.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'189.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.69%  "
$ws.Range("E7").Value = "  -2.93%  "
$ws.Range("D8").Value = "3.416.35"
$ws.Range("E8").Value = "  -5.38%  "
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("E10").Value = "  -5.63%  "
$ws.Range("E11").Value = "  -5.64%  "
$ws.Range("D12").Value = "'50.65"
$ws.Range("D12").Style = "Normal"
$ws.Range("E13").Value = "  -7.52%  "
$ws.Range("E14").Value = "  -5.77%  "
$ws.Range("D15").Value = "3.982.74"
$ws.Range("E15").Value = "  -5.15%  "
$ws.Range("D16").Value = "'626.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.37%  "
$ws.Range("D17").Value = "68.446.36"
$ws.Range("E17").Value = "  -2.90%  "
$ws.Range("D18").Value = "3.436.76"
$ws.Range("E18").Value = "  -5.59%  "
$ws.Range("E19").Value = "  -2.55%  "
$ws.Range("D20").Value = "'12.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.76%  "
$ws.Range("D21").Value = "'17.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.63%  "
$ws.Range("D22").Value = "'0.931"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.74%  "
$ws.Range("E23").Value = "  -2.65%  "
$ws.Range("E24").Value = "  +1.08%  "
$ws.Range("D25").Value = "'98.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.15%  "
$ws.Range("D26").Value = "'4.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -8.40%  "
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").Value = "'6.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.87%  "
$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").Value = "'2.81"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.72%  "
$ws.Range("D29").Value = "'9.75"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.85%  "
$ws.Range("E30").Value = "  -5.79%  "
$ws.Range("D31").Value = "'31.95"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.50%  "
$ws.Range("D32").Value = "'4.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -11.66%  "
$ws.Range("D33").Value = "'6.61"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.27%  "
$ws.Range("D34").Value = "'11.45"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.65%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "'60.49"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.29%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "'0.107"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.42%  "
$ws.Range("D38").Value = "3.644.68"
$ws.Range("E38").Value = "  -7.79%  "
$ws.Range("D39").Value = "0.0₃0769"
$ws.Range("E39").Value = "  -13.20%  "
$ws.Range("D40").Value = "'486.47"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.65%  "
$ws.Range("E41").Value = "  -2.75%  "
$ws.Range("B42").Value = "CoreDAO"
$ws.Range("C42").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D42").Value = "'3.60"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +75.48%  "
$ws.Range("B43").Value = "Fetch.AI"
$ws.Range("C43").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D43").Value = "'2.85"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.11%  "
$ws.Range("E44").Value = "  -6.36%  "
$ws.Range("E45").Value = "  -3.85%  "
$ws.Range("D46").Value = "'33.85"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.53%  "
$ws.Range("E47").Value = "  -6.20%  "
$ws.Range("D48").Value = "'3.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.37%  "
$ws.Range("D49").Value = "'2.77"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.32%  "
$ws.Range("E50").Value = "  -5.15%  "
$ws.Range("E51").Value = "  -0.18%  "
